# Generate Report for Archive
# - Localization status moves from "Ready for handoff" to "In Translation"
#   for the zh-cn / de-de targets (summarised on the Overview sheet and
#   detailed on each language sheet).
# - The Status column narrows to fit the shorter status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Update every cell that currently shows the old status ("Ready for
# handoff") so the shared string used across the sheets is replaced
# everywhere: the Overview summary (zh-cn/de-de columns) and each
# language sheet's own Status column.
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# Narrower status text means the Status column (and its Overview summary
# columns) no longer need to be as wide.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
